$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 3499.75
$ws.Range("J62").Value = 4000
$ws.Range("L62").Value = 4000
$ws.Range("N62").Value = -5248

# Row 65
$ws.Range("H65").Value = 3499.75
$ws.Range("J65").Value = 4000
$ws.Range("L65").Value = 20000
$ws.Range("N65").Value = -26240

# Row 76
$ws.Range("H76").Value = 7749.5
$ws.Range("J76").Value = 7999
$ws.Range("L76").Value = 7999
$ws.Range("N76").Value = -8629

# Row 79
$ws.Range("H79").Value = 7749.5
$ws.Range("J79").Value = 7999
$ws.Range("L79").Value = 7999
$ws.Range("N79").Value = -10183

# Row 112
$ws.Range("H112").Value = 2005

# Row 125
$ws.Range("H125").Value = 3356.111
$ws.Range("I125").Value = 1579.3334
$ws.Range("J125").Value = 4244.5
$ws.Range("K125").Value = 14214.0006
$ws.Range("L125").Value = 38200.5
$ws.Range("M125").Value = -11754.0006
$ws.Range("N125").Value = -43120.5

# Row 132
$ws.Range("H132").Value = 882.4091
$ws.Range("I132").Value = 837.29266
$ws.Range("K132").Value = 2511.87798
$ws.Range("M132").Value = 18.12202000000025

# Row 137
$ws.Range("H137").Value = 5014.684
$ws.Range("I137").Value = 3101.111
$ws.Range("K137").Value = 9303.332999999999
$ws.Range("M137").Value = -6753.332999999999

# Row 138
$ws.Range("H138").Value = 6419.4287
$ws.Range("J138").Value = 6640.4
$ws.Range("L138").Value = 19921.2
$ws.Range("N138").Value = -30201.2

$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 2784.35
$ws.Range("I74").Value = 2729
$ws.Range("J74").Value = 3098
$ws.Range("K74").Value = 2729
$ws.Range("L74").Value = 3098
$ws.Range("M74").Value = -1855
$ws.Range("N74").Value = -4846

# Row 77
$ws.Range("H77").Value = 2784.35
$ws.Range("I77").Value = 2729
$ws.Range("J77").Value = 3098
$ws.Range("K77").Value = 13645
$ws.Range("L77").Value = 15490
$ws.Range("M77").Value = -9277
$ws.Range("N77").Value = -24226

# Row 122
$ws.Range("H122").Value = 4004.55
$ws.Range("I122").Value = 2880.6875
$ws.Range("K122").Value = 8642.0625
$ws.Range("M122").Value = -6192.0625

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 1854.1666
$ws.Range("I20").Value = 1583.2222
$ws.Range("J20").Value = 2667
$ws.Range("K20").Value = 1583.2222
$ws.Range("L20").Value = 2667
$ws.Range("M20").Value = -1336.2222
$ws.Range("N20").Value = -3161

# Row 86
$ws.Range("H86").Value = 3429.8333
$ws.Range("I86").Value = 3515.8
$ws.Range("K86").Value = 3515.8
$ws.Range("M86").Value = -2392.8

# Row 89
$ws.Range("H89").Value = 3429.8333
$ws.Range("I89").Value = 3515.8
$ws.Range("K89").Value = 17579
$ws.Range("M89").Value = -11963

# Row 105
$ws.Range("H105").Value = 3963.2727
$ws.Range("I105").Value = 2371.4285
$ws.Range("K105").Value = 2371.4285
$ws.Range("M105").Value = -624.4285

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1053.3334
$ws.Range("I16").Value = 1200
$ws.Range("K16").Value = 1200
$ws.Range("M16").Value = -913

# Row 31
$ws.Range("H31").Value = 4195.4688
$ws.Range("I31").Value = 2876.4285
$ws.Range("K31").Value = 2876.4285
$ws.Range("M31").Value = -2581.4285

# Row 34
$ws.Range("H34").Value = 4195.4688
$ws.Range("I34").Value = 2876.4285
$ws.Range("K34").Value = 2876.4285
$ws.Range("M34").Value = -2674.4285

# Row 113
$ws.Range("H113").Value = 1053.3334
$ws.Range("I113").Value = 1200
$ws.Range("K113").Value = 1200
$ws.Range("M113").Value = 970

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 5659.857
$ws.Range("I5").Value = 587.6875
$ws.Range("K5").Value = 1763.0625
$ws.Range("M5").Value = -1651.0625

# Row 11
$ws.Range("H11").Value = 748.6667
$ws.Range("I11").Value = 664.8889
$ws.Range("K11").Value = 1994.6667
$ws.Range("M11").Value = -1854.6667

# Row 56
$ws.Range("H56").Value = 7715
$ws.Range("I56").Value = 7715
$ws.Range("K56").Value = 7715
$ws.Range("M56").Value = -7185

# Row 102
$ws.Range("H102").Value = 3963
$ws.Range("I102").Value = 426
$ws.Range("K102").Value = 1278
$ws.Range("M102").Value = 1156

# Row 113
$ws.Range("H113").Value = 2075.0625
$ws.Range("J113").Value = 2133.3845
$ws.Range("L113").Value = 6400.1535
$ws.Range("N113").Value = -10740.1535

# Row 122
$ws.Range("H122").Value = 1843.7142
$ws.Range("J122").Value = 2014.75
$ws.Range("L122").Value = 18132.75
$ws.Range("N122").Value = -23032.75

# Row 132
$ws.Range("H132").Value = 3288.8333
$ws.Range("I132").Value = 1384.3077
$ws.Range("J132").Value = 8240.6
$ws.Range("K132").Value = 12458.7693
$ws.Range("L132").Value = 74165.40000000001
$ws.Range("M132").Value = -9928.7693
$ws.Range("N132").Value = -79225.40000000001

# Row 135
$ws.Range("H135").Value = 5659.857
$ws.Range("I135").Value = 587.6875
$ws.Range("K135").Value = 5289.1875
$ws.Range("M135").Value = -2754.1875

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 25728918
$ws.Range("I80").Value = 126066.336
$ws.Range("K80").Value = 126066.336
$ws.Range("M80").Value = -125068.336

# Row 83
$ws.Range("H83").Value = 25728918
$ws.Range("I83").Value = 126066.336
$ws.Range("K83").Value = 630331.6799999999
$ws.Range("M83").Value = -625339.6799999999

# Row 122
$ws.Range("H122").Value = 3323.75
$ws.Range("I122").Value = 2244.6155
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 6733.8465
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = -4283.8465
$ws.Range("N122").Value = -28900

# Row 132
$ws.Range("H132").Value = 5049.2964
$ws.Range("I132").Value = 5133.24
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 15399.72
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -12869.72
$ws.Range("N132").Value = -17060

# Row 134
$ws.Range("H134").Value = 85387.664
$ws.Range("J134").Value = 85387.664
$ws.Range("L134").Value = 256162.992
$ws.Range("N134").Value = -261232.992

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1357
$ws.Range("I22").Value = 1357
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1357
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -1062
$ws.Range("N22").ClearContents()

# Row 27
$ws.Range("H27").Value = 1357
$ws.Range("I27").Value = 1357
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 1357
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -1250
$ws.Range("N27").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 2073.4167
$ws.Range("I126").Value = 2073.4167
$ws.Range("K126").Value = 6220.250100000001
$ws.Range("M126").Value = -3750.250100000001

# Row 136
$ws.Range("H136").Value = 3755.8684
$ws.Range("I136").Value = 2999.6
$ws.Range("J136").Value = 5210.231
$ws.Range("K136").Value = 8998.799999999999
$ws.Range("L136").Value = 15630.693
$ws.Range("M136").Value = -6448.799999999999
$ws.Range("N136").Value = -20730.693
